# "last minor changes to report"
#
# The sheet had four "summary" rows (8, 18, 28, 38) that carried
# leftover formula results from an earlier version of the model.
# The author cleared the stray numbers out of rows 8/18/28 (keeping the
# cell formatting that was already there) and removed row 38 entirely,
# which was the very last row of the sheet. Finally the active
# selection was left on the now-deleted-row's former location
# (A38:AZ38) instead of the old G13 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8, 18 and 28: wipe the formula/value content but leave the
# existing cell formatting (style) in place - same as selecting the
# row and pressing Delete / "Clear Contents".
$ws.Range("A8:AZ8").ClearContents()
$ws.Range("A18:AZ18").ClearContents()
$ws.Range("A28:AZ28").ClearContents()

# Row 38 (the old last row of the used range) is removed outright,
# shifting the sheet's dimension from AZ38 up to AZ37.
$ws.Rows.Item(38).Delete()

# Leave the selection where the removed row used to be, matching the
# saved workbook's cursor position.
$ws.Range("A38:AZ38").Select()
